$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = "June 28, 2023"
$ws.Cells.Item(2, 4).Value = 82862
$ws.Cells.Item(2, 6).Value = 61798
$ws.Cells.Item(2, 9).Value = 21064

$ws.Cells.Item(3, 1).Value = "June 28, 2023"
$ws.Cells.Item(3, 4).Value = 4746
$ws.Cells.Item(3, 9).Value = 317

$ws.Cells.Item(4, 1).Value = "June 28, 2023"

$ws.Cells.Item(5, 1).Value = "June 28, 2023"
$ws.Cells.Item(5, 4).Value = 489751
$ws.Cells.Item(5, 6).Value = 399366
$ws.Cells.Item(5, 7).Value = 101
$ws.Cells.Item(5, 9).Value = 90385

$ws.Cells.Item(6, 1).Value = "June 28, 2023"
$ws.Cells.Item(6, 4).Value = 5752049
$ws.Cells.Item(6, 6).Value = 4755313
$ws.Cells.Item(6, 9).Value = 996736

$ws.Cells.Item(7, 1).Value = "June 28, 2023"
$ws.Cells.Item(7, 4).Value = 233570
$ws.Cells.Item(7, 6).Value = 43177
$ws.Cells.Item(7, 7).Value = 96
$ws.Cells.Item(7, 9).Value = 190393

$ws.Cells.Item(8, 1).Value = "June 28, 2023"
$ws.Cells.Item(8, 4).Value = 381983
$ws.Cells.Item(8, 6).Value = 322339
$ws.Cells.Item(8, 9).Value = 59644

$ws.Cells.Item(9, 1).Value = "June 28, 2023"
$ws.Cells.Item(9, 4).Value = 3187439
$ws.Cells.Item(9, 6).Value = 2412836
$ws.Cells.Item(9, 9).Value = 774603

$ws.Cells.Item(10, 1).Value = "June 28, 2023"
$ws.Cells.Item(10, 6).Value = 42793
$ws.Cells.Item(10, 7).Value = 92
$ws.Cells.Item(10, 9).Value = 12614

$ws.Cells.Item(11, 1).Value = "June 28, 2023"
$ws.Cells.Item(11, 4).Value = 30392
$ws.Cells.Item(11, 6).Value = 24997
$ws.Cells.Item(11, 7).Value = 71
$ws.Cells.Item(11, 9).Value = 5395

$ws.Cells.Item(12, 1).Value = "June 28, 2023"
$ws.Cells.Item(12, 4).Value = 2648232
$ws.Cells.Item(12, 6).Value = 2011630
$ws.Cells.Item(12, 9).Value = 636602

$ws.Cells.Item(13, 1).Value = "June 28, 2023"
$ws.Cells.Item(13, 4).Value = 1190773
$ws.Cells.Item(13, 6).Value = 1047294
$ws.Cells.Item(13, 9).Value = 143479

$ws.Cells.Item(14, 1).Value = "June 28, 2023"
$ws.Cells.Item(14, 4).Value = 334964
$ws.Cells.Item(14, 6).Value = 258155
$ws.Cells.Item(14, 7).Value = 113
$ws.Cells.Item(14, 9).Value = 76809

$ws.Cells.Item(15, 1).Value = "June 28, 2023"
$ws.Cells.Item(15, 4).Value = 171488
$ws.Cells.Item(15, 6).Value = 116111
$ws.Cells.Item(15, 7).Value = 189
$ws.Cells.Item(15, 9).Value = 55377

$ws.Cells.Item(16, 1).Value = "June 28, 2023"
$ws.Cells.Item(16, 4).Value = 106931
$ws.Cells.Item(16, 6).Value = 77293
$ws.Cells.Item(16, 7).Value = 152
$ws.Cells.Item(16, 9).Value = 29638

$ws.Cells.Item(17, 1).Value = "June 28, 2023"
$ws.Cells.Item(17, 4).Value = 428
$ws.Cells.Item(17, 9).Value = 253

$ws.Cells.Item(18, 1).Value = "June 28, 2023"
$ws.Cells.Item(18, 4).Value = 601378
$ws.Cells.Item(18, 6).Value = 411965
$ws.Cells.Item(18, 7).Value = 53
$ws.Cells.Item(18, 9).Value = 189413

$ws.Cells.Item(19, 1).Value = "June 28, 2023"
$ws.Cells.Item(19, 6).Value = 3991280
$ws.Cells.Item(19, 7).Value = 14
$ws.Cells.Item(19, 9).Value = 277555

$ws.Cells.Item(20, 1).Value = "June 28, 2023"
$ws.Cells.Item(20, 7).Value = 158

$ws.Cells.Item(21, 1).Value = "June 28, 2023"

$ws.Cells.Item(22, 1).Value = "June 28, 2023"
$ws.Cells.Item(22, 4).Value = 382831
$ws.Cells.Item(22, 7).Value = 59
$ws.Cells.Item(22, 9).Value = 379702

$ws.Cells.Item(23, 1).Value = "June 28, 2023"
$ws.Cells.Item(23, 4).Value = 81595
$ws.Cells.Item(23, 7).Value = 192
$ws.Cells.Item(23, 9).Value = 36372

$ws.Cells.Item(24, 1).Value = "June 28, 2023"
$ws.Cells.Item(24, 4).Value = 1031662
$ws.Cells.Item(24, 6).Value = 877815
$ws.Cells.Item(24, 9).Value = 153847

$ws.Cells.Item(25, 1).Value = "June 28, 2023"
$ws.Cells.Item(25, 3).Value = 3.7
$ws.Cells.Item(25, 4).Value = 1455
$ws.Cells.Item(25, 6).Value = 303
$ws.Cells.Item(25, 8).Value = 3.8

$ws.Cells.Item(26, 1).Value = "June 28, 2023"
$ws.Cells.Item(26, 4).Value = 4578
$ws.Cells.Item(26, 6).Value = 1625
$ws.Cells.Item(26, 9).Value = 2953

$ws.Cells.Item(27, 1).Value = "June 28, 2023"
$ws.Cells.Item(27, 4).Value = 475
$ws.Cells.Item(27, 9).Value = 461

$ws.Cells.Item(28, 1).Value = "June 28, 2023"

$ws.Cells.Item(29, 1).Value = "June 28, 2023"
$ws.Cells.Item(29, 4).Value = 1847
$ws.Cells.Item(29, 9).Value = 1540

$ws.Cells.Item(30, 1).Value = "June 28, 2023"

$ws.Cells.Item(31, 1).Value = "June 28, 2023"

$ws.Cells.Item(32, 1).Value = "June 28, 2023"
$ws.Cells.Item(32, 4).Value = 2601
$ws.Cells.Item(32, 6).Value = 1377

$ws.Cells.Item(33, 1).Value = "June 28, 2023"
$ws.Cells.Item(33, 4).Value = 48931
$ws.Cells.Item(33, 6).Value = 28397
$ws.Cells.Item(33, 9).Value = 20534

$ws.Cells.Item(34, 1).Value = "June 28, 2023"
$ws.Cells.Item(34, 3).Value = 2.45
$ws.Cells.Item(34, 4).Value = 32
$ws.Cells.Item(34, 5).Value = 3.1
$ws.Cells.Item(34, 6).Value = 9

$ws.Cells.Item(35, 1).Value = "June 28, 2023"
$ws.Cells.Item(35, 4).Value = 15909
$ws.Cells.Item(35, 6).Value = 13446

$ws.Cells.Item(36, 1).Value = "June 28, 2023"

$ws.Cells.Item(37, 1).Value = "June 28, 2023"
$ws.Cells.Item(37, 4).Value = 4184
$ws.Cells.Item(37, 9).Value = 884

$ws.Cells.Item(38, 1).Value = "June 28, 2023"
